# Update on 12/10/2025 at 8:58pm
$wb = $excel.ActiveWorkbook

# --- "Table" sheet: refresh the FY26 (YTD) row with the latest reported data ---
$tableWs = $wb.Worksheets.Item("Table")

$tableWs.Range("B9").Value = 7
$tableWs.Range("C9").Value = 469889
$tableWs.Range("D9").Value = 67873
$tableWs.Range("E9").Value = 2108388
$tableWs.Range("F9").Value = 2646150
$tableWs.Range("G9").Value = 378021.428571428
$tableWs.Range("H9").Value = 377032
$tableWs.Range("I9").Value = 0.0011790586541767299

# Unhide column B (Months Reported) on the Table sheet
$tableWs.Columns.Item(2).Hidden = $false

# Update selection / view on Table sheet
$tableWs.Activate()
$tableWs.Range("C10").Select()
$excel.ActiveWindow.ScrollRow = 1

# --- "projection_data" sheet (hidden): update saved view/selection state ---
$projWs = $wb.Worksheets.Item("projection_data")
$projWs.Visible = $true
$projWs.Activate()
$projWs.Range("D13").Select()
$excel.ActiveWindow.ScrollRow = 6
$projWs.Visible = $false

# Restore the originally active sheet ("Table", tabSelected)
$tableWs.Activate()

$wb.Application.Calculate()
